# Auto-generated edit script: updates market-price derived columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to refreshed values from a
# scheduled market-data pull. Pure value writes, no formulas/styles touched.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1062.25
$ws.Range("J43").Value = 1999
$ws.Range("L43").Value = 1999
$ws.Range("N43").Value = -2137
$ws.Range("H53").Value = 359.18182
$ws.Range("I53").Value = 316.55554
$ws.Range("K53").Value = 316.55554
$ws.Range("M53").Value = 320.44446
$ws.Range("H76").Value = 4458.25
$ws.Range("I76").Value = 4166.5
$ws.Range("J76").Value = 4750
$ws.Range("K76").Value = 4166.5
$ws.Range("L76").Value = 4750
$ws.Range("M76").Value = -3851.5
$ws.Range("N76").Value = -5380
$ws.Range("H79").Value = 4458.25
$ws.Range("I79").Value = 4166.5
$ws.Range("J79").Value = 4750
$ws.Range("K79").Value = 4166.5
$ws.Range("L79").Value = 4750
$ws.Range("M79").Value = -3074.5
$ws.Range("N79").Value = -6934
$ws.Range("H125").Value = 4480.273
$ws.Range("I125").Value = 3108
$ws.Range("J125").Value = 8139.6665
$ws.Range("K125").Value = 27972
$ws.Range("L125").Value = 73256.9985
$ws.Range("M125").Value = -25512
$ws.Range("N125").Value = -78176.9985
$ws.Range("H138").Value = 2499.2942
$ws.Range("I138").Value = 497.6
$ws.Range("K138").Value = 1492.8
$ws.Range("M138").Value = 3647.2
$ws.Range("H141").Value = 4239.4546
$ws.Range("I141").Value = 4239.4546
$ws.Range("K141").Value = 12718.3638
$ws.Range("M141").Value = -7538.363799999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 849
$ws.Range("I4").Value = 998
$ws.Range("K4").Value = 998
$ws.Range("M4").Value = -882
$ws.Range("H32").Value = 8429.777
$ws.Range("I32").Value = 3337.4119
$ws.Range("K32").Value = 3337.4119
$ws.Range("M32").Value = -3050.4119
$ws.Range("H61").Value = 2337.5454
$ws.Range("I61").Value = 2337.5454
$ws.Range("K61").Value = 2337.5454
$ws.Range("M61").Value = -2125.5454
$ws.Range("H132").Value = 1916.2
$ws.Range("I132").Value = 1911.7894
$ws.Range("K132").Value = 5735.3682
$ws.Range("M132").Value = -3205.3682
$ws.Range("H136").Value = 2337.5454
$ws.Range("I136").Value = 2337.5454
$ws.Range("K136").Value = 7012.6362
$ws.Range("M136").Value = -4462.6362

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2727.6843
$ws.Range("I94").Value = 2754.4707
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 2754.4707
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -2303.4707
$ws.Range("N94").Value = -3402
$ws.Range("H134").Value = 3104.6956
$ws.Range("I134").Value = 3104.6956
$ws.Range("K134").Value = 9314.086800000001
$ws.Range("M134").Value = -6779.086800000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2869.5757
$ws.Range("I31").Value = 2349.85
$ws.Range("K31").Value = 2349.85
$ws.Range("M31").Value = -2054.85
$ws.Range("H34").Value = 2869.5757
$ws.Range("I34").Value = 2349.85
$ws.Range("K34").Value = 2349.85
$ws.Range("M34").Value = -2147.85
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H99").Value = 6114.143
$ws.Range("I99").Value = 6059.8
$ws.Range("J99").Value = 6250
$ws.Range("K99").Value = 6059.8
$ws.Range("L99").Value = 6250
$ws.Range("M99").Value = -4561.8
$ws.Range("N99").Value = -9246
$ws.Range("H122").Value = 532.2308
$ws.Range("I122").Value = 509.91666
$ws.Range("K122").Value = 1529.74998
$ws.Range("M122").Value = 920.2500199999999
$ws.Range("H126").Value = 6114.143
$ws.Range("I126").Value = 6059.8
$ws.Range("J126").Value = 6250
$ws.Range("K126").Value = 18179.4
$ws.Range("L126").Value = 18750
$ws.Range("M126").Value = -15709.4
$ws.Range("N126").Value = -23690
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 226.33333
$ws.Range("J13").Value = 19.666666
$ws.Range("L13").Value = 58.999998
$ws.Range("N13").Value = -394.999998
$ws.Range("H34").Value = 435.81818
$ws.Range("J34").Value = 500
$ws.Range("L34").Value = 1500
$ws.Range("N34").Value = -1668
$ws.Range("H39").Value = 4850
$ws.Range("J39").Value = 4850
$ws.Range("L39").Value = 14550
$ws.Range("N39").Value = -15138
$ws.Range("H55").Value = 3340.9092
$ws.Range("J55").Value = 4678.5713
$ws.Range("L55").Value = 14035.7139
$ws.Range("N55").Value = -14389.7139

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 9999
$ws.Range("I132").Value = 9999
$ws.Range("K132").Value = 29997
$ws.Range("M132").Value = -27467

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1999.6
$ws.Range("I4").Value = 1999.6
$ws.Range("K4").Value = 1999.6
$ws.Range("M4").Value = -1886.6
$ws.Range("H28").Value = 1999.6
$ws.Range("I28").Value = 1999.6
$ws.Range("K28").Value = 1999.6
$ws.Range("M28").Value = -1767.6
$ws.Range("H37").Value = 1999.6
$ws.Range("I37").Value = 1999.6
$ws.Range("K37").Value = 1999.6
$ws.Range("M37").Value = -1892.6
$ws.Range("H40").Value = 3599.6
$ws.Range("I40").Value = 4666.3335
$ws.Range("J40").Value = 1999.5
$ws.Range("K40").Value = 4666.3335
$ws.Range("L40").Value = 1999.5
$ws.Range("M40").Value = -4530.3335
$ws.Range("N40").Value = -2271.5
$ws.Range("H43").Value = 12000
$ws.Range("I43").Value = 12000
$ws.Range("K43").Value = 12000
$ws.Range("M43").Value = -11807
$ws.Range("H132").Value = 6038.0835
$ws.Range("I132").Value = 2559.8
$ws.Range("J132").Value = 8522.571
$ws.Range("K132").Value = 7679.400000000001
$ws.Range("L132").Value = 25567.713
$ws.Range("M132").Value = -5149.400000000001
$ws.Range("N132").Value = -30627.713

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1973.3636
$ws.Range("I132").Value = 1973.3636
$ws.Range("K132").Value = 5920.0908
$ws.Range("M132").Value = -3390.0908
